# Updates the cryptos list worksheet, replacing price and volume change
# values per the latest scrape, and swapping the Toncoin / LidoDAOToken
# rows (25/26) to reflect the new ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($sheet, [string]$addr, [string]$text)
    $cell = $sheet.Range($addr)
    # Force the cell to be stored as text (matches the source data, which
    # is all inline/shared string content) so values like "1.000" or
    # "24.949.48" are not reinterpreted as numbers/dates.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    # Restore the default (unstyled) cell style so no stray formatting
    # is introduced compared to the original workbook.
    $cell.Style = "Normal"
}

Set-CellText $ws "D2" "24.949.48"
Set-CellText $ws "E2" "  +0.21%  "
Set-CellText $ws "D3" "1.708.23"
Set-CellText $ws "E3" "  -0.19%  "
Set-CellText $ws "D4" "1.000"
Set-CellText $ws "E4" "  -0.46%  "
Set-CellText $ws "D5" "318.05"
Set-CellText $ws "E5" "  +0.13%  "
Set-CellText $ws "D6" "1.001"
Set-CellText $ws "E6" "  -0.20%  "
Set-CellText $ws "D7" "0.4035"
Set-CellText $ws "E7" "  +1.80%  "
Set-CellText $ws "D8" "0.4074"
Set-CellText $ws "E8" "  -0.89%  "
Set-CellText $ws "D9" "1.481"
Set-CellText $ws "E9" "  -2.87%  "
Set-CellText $ws "D10" "53.85"
Set-CellText $ws "E10" "  +0.51%  "
Set-CellText $ws "D11" "0.9999"
Set-CellText $ws "E11" "  -0.53%  "
Set-CellText $ws "D12" "0.08833"
Set-CellText $ws "E12" "  -1.12%  "
Set-CellText $ws "D13" "26.33"
Set-CellText $ws "E13" "  +7.06%  "
Set-CellText $ws "D14" "7.523"
Set-CellText $ws "E14" "  -2.51%  "
Set-CellText $ws "D15" "8.137"
Set-CellText $ws "E15" "  -0.30%  "
Set-CellText $ws "D16" "0.00001360"
Set-CellText $ws "E16" "  -1.04%  "
Set-CellText $ws "D17" "1.738.08"
Set-CellText $ws "E17" "  +2.60%  "
Set-CellText $ws "D18" "96.97"
Set-CellText $ws "E18" "  -3.42%  "
Set-CellText $ws "D19" "0.07158"
Set-CellText $ws "E19" "  +0.24%  "
Set-CellText $ws "D20" "21.21"
Set-CellText $ws "E20" "  +5.36%  "
Set-CellText $ws "D21" "7.300"
Set-CellText $ws "E21" "  -2.58%  "
Set-CellText $ws "E22" "  -0.51%  "
Set-CellText $ws "D23" "14.41"
Set-CellText $ws "E23" "  -0.60%  "
Set-CellText $ws "D24" "24.938.87"
Set-CellText $ws "E24" "  +0.12%  "
Set-CellText $ws "B25" "Toncoin"
Set-CellText $ws "C25" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-CellText $ws "D25" "2.334"
Set-CellText $ws "E25" "  +0.02%  "
Set-CellText $ws "B26" "LidoDAOToken"
Set-CellText $ws "C26" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-CellText $ws "D26" "2.920"
Set-CellText $ws "E26" "  -6.57%  "
Set-CellText $ws "D27" "23.30"
Set-CellText $ws "E27" "  +1.05%  "
Set-CellText $ws "D28" "6.245"
Set-CellText $ws "E28" "  +19.20%  "
Set-CellText $ws "D29" "167.38"
Set-CellText $ws "E29" "  +1.20%  "
Set-CellText $ws "D30" "146.39"
Set-CellText $ws "E30" "  +4.00%  "
Set-CellText $ws "D31" "8.393"
Set-CellText $ws "E31" "  -10.44%  "
Set-CellText $ws "D32" "1.923.57"
Set-CellText $ws "E32" "  +2.23%  "
Set-CellText $ws "D33" "2.233"
Set-CellText $ws "E33" "  +13.44%  "
Set-CellText $ws "D34" "0.08911"
Set-CellText $ws "E34" "  -1.50%  "
Set-CellText $ws "D35" "0.03213"
Set-CellText $ws "E35" "  +6.39%  "
Set-CellText $ws "D36" "7.270"
Set-CellText $ws "E36" "  -7.66%  "
Set-CellText $ws "D37" "1.029"
Set-CellText $ws "E37" "  -5.47%  "
Set-CellText $ws "D38" "0.2862"
Set-CellText $ws "E38" "  +1.90%  "
Set-CellText $ws "D39" "0.8503"
Set-CellText $ws "E39" "  +5.34%  "
Set-CellText $ws "D40" "10.90"
Set-CellText $ws "E40" "  -1.70%  "
Set-CellText $ws "D41" "0.09357"
Set-CellText $ws "E41" "  +0.52%  "
Set-CellText $ws "D42" "14.20"
Set-CellText $ws "E42" "  -2.61%  "
Set-CellText $ws "E43" "  -1.21%  "
Set-CellText $ws "D44" "17.54"
Set-CellText $ws "E44" "  +5.06%  "
Set-CellText $ws "D45" "2.723"
Set-CellText $ws "E45" "  +3.12%  "
Set-CellText $ws "D46" "0.7459"
Set-CellText $ws "E46" "  +1.24%  "
Set-CellText $ws "D47" "4.255"
Set-CellText $ws "E47" "  -0.33%  "
Set-CellText $ws "D48" "1.400"
Set-CellText $ws "E48" "  +3.85%  "
Set-CellText $ws "D49" "1.000"
Set-CellText $ws "E49" "  -0.11%  "
Set-CellText $ws "D50" "142.47"
Set-CellText $ws "E50" "  +1.09%  "
Set-CellText $ws "D51" "0.08378"
Set-CellText $ws "E51" "  +3.56%  "
